$d = $word.ActiveDocument

$replacements = @(
    @("91×83=", "23×64="),
    @("76×14=", "83×53="),
    @("75×13=", "69×28="),
    @("92×76=", "26×21="),
    @("93×53=", "82×42="),
    @("34×59=", "78×46="),
    @("96×44=", "67×47="),
    @("19×54=", "94×59="),
    @("47×13=", "59×31="),
    @("85×13=", "21×71="),
    @("75×25=", "52×65="),
    @("26×62=", "74×55="),
    @("71×92=", "32×42="),
    @("14×86=", "58×86="),
    @("82×48=", "49×16="),
    @("90×14=", "76×16="),
    @("95×90=", "61×61="),
    @("39×43=", "48×96="),
    @("56×14=", "40×91="),
    @("63×47=", "21×80="),
    @("60×39=", "80×47="),
    @("42×24=", "35×77="),
    @("78×83=", "78×51="),
    @("97×70=", "77×65="),
    @("90×28=", "99×47=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
